# "Lista de Restrições" – atualização da tabela de alunos.
#
# O aluno "Brando de Oliveira Veridiano" (RA 1901003) saiu do grupo; sua
# linha inteira é removida da tabela de integrantes. A linha seguinte
# ("Guilherme Alves dos Santos", RA 1900785, com seu e-mail e celular)
# passa a ser a primeira linha de dados da tabela.

$d = $word.ActiveDocument

# A tabela de integrantes é a primeira tabela do documento.
$t = $d.Tables.Item(1)

# Linha 1 = cabeçalho (Aluno/RA/E-mail/Celular); linha 2 = o aluno que saiu.
$row = $t.Rows.Item(2)
Write-Host "Removendo linha:" $row.Range.Text
$row.Delete()

Write-Host "Linhas restantes na tabela:" $t.Rows.Count
